# Update countries & provincias Spain
#
# 1) Refresh the "last updated" timestamp on the report header.
# 2) Re-order a handful of low-case-count countries in the list (their
#    whole data rows move together with the name), which is how the
#    author's shared-string reshuffle manifests in the sheet.
# 3) Refresh case counts for Rusia, Oman, Moldavia and Eslovaquia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 10:10"

# --- helper: swap two full data rows (columns A:H) ----------------------
function Swap-Rows($r1, $r2) {
    for ($c = 1; $c -le 8; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# --- 2) Re-order countries ----------------------------------------------
# Belice / Nueva Caledonia swap places
Swap-Rows 200 201
# Papua Nueva Guinea / Islas Virgenes Britanicas swap places
Swap-Rows 213 214
# San Bartolome / Bonaire, San Eustaquio y Saba swap places
Swap-Rows 215 216

# --- 3) Updated case figures ---------------------------------------------
# Rusia (row 6)
$ws.Cells.Item(6, 2).Value = 379051
$ws.Cells.Item(6, 3).Value = 8371
$ws.Cells.Item(6, 4).Value = 150993
$ws.Cells.Item(6, 5).Value = 223916
$ws.Cells.Item(6, 7).Value = 174
$ws.Cells.Item(6, 8).Value = 4142

# Oman (row 59)
$ws.Cells.Item(59, 5).Value = 6156
$ws.Cells.Item(59, 7).Value = 1
$ws.Cells.Item(59, 8).Value = 40

# Moldavia (row 64)
$ws.Cells.Item(64, 4).Value = 4123
$ws.Cells.Item(64, 5).Value = 3138
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = 276

# Eslovaquia (row 98)
$ws.Cells.Item(98, 2).Value = 1520
$ws.Cells.Item(98, 3).Value = 5
$ws.Cells.Item(98, 4).Value = 1332
